$d = $word.ActiveDocument

$d.Content.Find.Execute("98×34=3332", $true, $false, $false, $false, $false, $true, 1, $false, "65×63=4095", 2)
$d.Content.Find.Execute("96×40=3840", $true, $false, $false, $false, $false, $true, 1, $false, "93×78=7254", 2)
$d.Content.Find.Execute("18×56=1008", $true, $false, $false, $false, $false, $true, 1, $false, "60×27=1620", 2)
$d.Content.Find.Execute("13×50=650", $true, $false, $false, $false, $false, $true, 1, $false, "62×18=1116", 2)
$d.Content.Find.Execute("49×23=1127", $true, $false, $false, $false, $false, $true, 1, $false, "64×11=704", 2)
$d.Content.Find.Execute("66×25=1650", $true, $false, $false, $false, $false, $true, 1, $false, "54×18=972", 2)
$d.Content.Find.Execute("55×15=825", $true, $false, $false, $false, $false, $true, 1, $false, "79×32=2528", 2)
$d.Content.Find.Execute("41×37=1517", $true, $false, $false, $false, $false, $true, 1, $false, "22×29=638", 2)
$d.Content.Find.Execute("53×13=689", $true, $false, $false, $false, $false, $true, 1, $false, "67×33=2211", 2)
$d.Content.Find.Execute("68×51=3468", $true, $false, $false, $false, $false, $true, 1, $false, "27×65=1755", 2)
$d.Content.Find.Execute("97×12=1164", $true, $false, $false, $false, $false, $true, 1, $false, "70×24=1680", 2)
$d.Content.Find.Execute("87×71=6177", $true, $false, $false, $false, $false, $true, 1, $false, "32×81=2592", 2)
$d.Content.Find.Execute("36×91=3276", $true, $false, $false, $false, $false, $true, 1, $false, "70×90=6300", 2)
$d.Content.Find.Execute("77×73=5621", $true, $false, $false, $false, $false, $true, 1, $false, "38×65=2470", 2)
$d.Content.Find.Execute("77×97=7469", $true, $false, $false, $false, $false, $true, 1, $false, "34×99=3366", 2)
$d.Content.Find.Execute("32×94=3008", $true, $false, $false, $false, $false, $true, 1, $false, "11×15=165", 2)
$d.Content.Find.Execute("31×69=2139", $true, $false, $false, $false, $false, $true, 1, $false, "39×76=2964", 2)
$d.Content.Find.Execute("42×42=1764", $true, $false, $false, $false, $false, $true, 1, $false, "97×53=5141", 2)
$d.Content.Find.Execute("15×50=750", $true, $false, $false, $false, $false, $true, 1, $false, "15×27=405", 2)
$d.Content.Find.Execute("54×51=2754", $true, $false, $false, $false, $false, $true, 1, $false, "69×67=4623", 2)
$d.Content.Find.Execute("48×78=3744", $true, $false, $false, $false, $false, $true, 1, $false, "43×12=516", 2)
$d.Content.Find.Execute("87×57=4959", $true, $false, $false, $false, $false, $true, 1, $false, "33×44=1452", 2)
$d.Content.Find.Execute("68×89=6052", $true, $false, $false, $false, $false, $true, 1, $false, "40×71=2840", 2)
$d.Content.Find.Execute("24×73=1752", $true, $false, $false, $false, $false, $true, 1, $false, "64×55=3520", 2)
$d.Content.Find.Execute("69×28=1932", $true, $false, $false, $false, $false, $true, 1, $false, "72×26=1872", 2)
